$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.019999999999999
$ws.Range("C2").Value = 1.037242964115984
$ws.Range("D2").Value = 1.045846838477866
$ws.Range("E2").Value = 1.036062008959118
$ws.Range("F2").Value = 1.056812275636661
$ws.Range("I2").Value = 1.043584320200308
$ws.Range("J2").Value = 1.042346897455935
$ws.Range("K2").Value = 1.048613931761185
$ws.Range("L2").Value = 1.038856826762764
$ws.Range("M2").Value = 1.059548960345595
$ws.Range("N2").Value = 1.043827149609183

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.038093349774712
$ws.Range("D3").Value = 1.046519585135164
$ws.Range("E3").Value = 1.03678250286227
$ws.Range("F3").Value = 1.057653666148145
$ws.Range("I3").Value = 1.043824413062562
$ws.Range("J3").Value = 1.042842079608715
$ws.Range("K3").Value = 1.049098562681563
$ws.Range("L3").Value = 1.039387080923744
$ws.Range("M3").Value = 1.060203993285993
$ws.Range("N3").Value = 1.044323034977418

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038644175181942
$ws.Range("D4").Value = 1.046955368983055
$ws.Range("E4").Value = 1.037249576455663
$ws.Range("F4").Value = 1.058198951808635
$ws.Range("I4").Value = 1.043978760668798
$ws.Range("J4").Value = 1.043162402411283
$ws.Range("K4").Value = 1.049411940765784
$ws.Range("L4").Value = 1.039730385204654
$ws.Range("M4").Value = 1.060628044169536
$ws.Range("N4").Value = 1.044643812675107

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.038875876846415
$ws.Range("D5").Value = 1.047138683729103
$ws.Range("E5").Value = 1.037446139725525
$ws.Range("F5").Value = 1.058428391440818
$ws.Range("I5").Value = 1.044043406149814
$ws.Range("J5").Value = 1.043297042576971
$ws.Range("K5").Value = 1.04954363322608
$ws.Range("L5").Value = 1.039874755524607
$ws.Range("M5").Value = 1.060806361634949
$ws.Range("N5").Value = 1.044778644045273

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.038914788426794
$ws.Range("D6").Value = 1.047169469553113
$ws.Range("E6").Value = 1.037479155564776
$ws.Range("F6").Value = 1.058466927128178
$ws.Range("I6").Value = 1.044054246181247
$ws.Range("J6").Value = 1.043319647836843
$ws.Range("K6").Value = 1.049565741901355
$ws.Range("L6").Value = 1.039898998548328
$ws.Range("M6").Value = 1.060836304601036
$ws.Range("N6").Value = 1.044801281407207

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.03864727066374
$ws.Range("D7").Value = 1.046957818008831
$ws.Range("E7").Value = 1.037252202138309
$ws.Range("F7").Value = 1.058202016802683
$ws.Range("I7").Value = 1.043979625417775
$ws.Range("J7").Value = 1.043164201573817
$ws.Range("K7").Value = 1.049413700651588
$ws.Range("L7").Value = 1.039732314112006
$ws.Range("M7").Value = 1.060630426675621
$ws.Range("N7").Value = 1.044645614392658

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.037530237171117
$ws.Range("D8").Value = 1.046074097638079
$ws.Range("E8").Value = 1.036305322908821
$ws.Range("F8").Value = 1.057096450420338
$ws.Range("I8").Value = 1.043665669043544
$ws.Range("J8").Value = 1.042514265128786
$ws.Range("K8").Value = 1.04877775786797
$ws.Range("L8").Value = 1.039035987588543
$ws.Range("M8").Value = 1.059770289267345
$ws.Range("N8").Value = 1.043994754963326

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03556630858921
$ws.Range("D9").Value = 1.044520559004271
$ws.Range("E9").Value = 1.034643509846727
$ws.Range("F9").Value = 1.055154892667639
$ws.Range("I9").Value = 1.043104749445816
$ws.Range("J9").Value = 1.041368335260888
$ws.Range("K9").Value = 1.047655589870731
$ws.Range("L9").Value = 1.037810519028627
$ws.Range("M9").Value = 1.058256225253859
$ws.Range("N9").Value = 1.042847197743576

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.034260092270597
$ws.Range("D10").Value = 1.043487450153435
$ws.Range("E10").Value = 1.033540250061773
$ws.Range("F10").Value = 1.053865065067088
$ws.Range("I10").Value = 1.042725681986404
$ws.Range("J10").Value = 1.040604012702413
$ws.Range("K10").Value = 1.04690650853476
$ws.Range("L10").Value = 1.03699466011456
$ws.Range("M10").Value = 1.057248025115485
$ws.Range("N10").Value = 1.042081789759398

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.033695233667058
$ws.Range("D11").Value = 1.043040737115542
$ws.Range("E11").Value = 1.033063643270504
$ws.Range("F11").Value = 1.053307657021933
$ws.Range("I11").Value = 1.04256033705398
$ws.Range("J11").Value = 1.040272980490282
$ws.Range("K11").Value = 1.046581934241339
$ws.Range("L11").Value = 1.03664166663216
$ws.Range("M11").Value = 1.056811762180315
$ws.Range("N11").Value = 1.041750287443557

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.033485533013863
$ws.Range("D12").Value = 1.042874904393918
$ws.Range("E12").Value = 1.032886779141842
$ws.Range("F12").Value = 1.053100777403712
$ws.Range("I12").Value = 1.042498740157908
$ws.Range("J12").Value = 1.040150010283843
$ws.Range("K12").Value = 1.046461341564643
$ws.Range("L12").Value = 1.036510592372273
$ws.Range("M12").Value = 1.056649760531699
$ws.Range("N12").Value = 1.04162714260532

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.033530509378153
$ws.Range("D13").Value = 1.042910471685631
$ws.Range("E13").Value = 1.032924709428034
$ws.Range("F13").Value = 1.053145146216695
$ws.Range("I13").Value = 1.042511961053095
$ws.Range("J13").Value = 1.040176388247864
$ws.Range("K13").Value = 1.046487210500316
$ws.Range("L13").Value = 1.036538706267777
$ws.Range("M13").Value = 1.056684508335605
$ws.Range("N13").Value = 1.041653558029076

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.033677897431845
$ws.Range("D14").Value = 1.043027027354569
$ws.Range("E14").Value = 1.033049020171514
$ws.Range("F14").Value = 1.053290552879549
$ws.Range("I14").Value = 1.04255524911137
$ws.Range("J14").Value = 1.040262815933599
$ws.Range("K14").Value = 1.046571966644781
$ws.Range("L14").Value = 1.036630831103142
$ws.Range("M14").Value = 1.056798370129311
$ws.Range("N14").Value = 1.041740108452037

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.033768723093939
$ws.Range("D15").Value = 1.043098853974123
$ws.Range("E15").Value = 1.033125634551797
$ws.Range("F15").Value = 1.053380164850083
$ws.Range("I15").Value = 1.04258189642835
$ws.Range("J15").Value = 1.040316065581503
$ws.Range("K15").Value = 1.046624183589401
$ws.Range("L15").Value = 1.03668759802429
$ws.Range("M15").Value = 1.056868530259219
$ws.Range("N15").Value = 1.04179343372055

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.034297595834683
$ws.Range("D16").Value = 1.043517110447778
$ws.Range("E16").Value = 1.033571904465969
$ws.Range("F16").Value = 1.053902081654276
$ws.Range("I16").Value = 1.042736630039135
$ws.Range("J16").Value = 1.040625980735185
$ws.Range("K16").Value = 1.046928044996975
$ws.Range("L16").Value = 1.037018093118764
$ws.Range("M16").Value = 1.05727698478874
$ws.Range("N16").Value = 1.042103788989297

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.034629543397046
$ws.Range("D17").Value = 1.0437796415672
$ws.Range("E17").Value = 1.033852136766541
$ws.Range("F17").Value = 1.054229761009464
$ws.Range("I17").Value = 1.042833367966103
$ws.Range("J17").Value = 1.040820362984481
$ws.Range("K17").Value = 1.047118592007721
$ws.Range("L17").Value = 1.037225479557633
$ws.Range("M17").Value = 1.057533277468885
$ws.Range("N17").Value = 1.042298447283684

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.034823234251302
$ws.Range("D18").Value = 1.043932832187373
$ws.Range("E18").Value = 1.034015698695757
$ws.Range("F18").Value = 1.054420996435587
$ws.Range("I18").Value = 1.042889677125979
$ws.Range("J18").Value = 1.040933735493965
$ws.Range("K18").Value = 1.047229713818765
$ws.Range("L18").Value = 1.037346471369581
$ws.Range("M18").Value = 1.057682796948762
$ws.Range("N18").Value = 1.042411980795135

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.034889289900367
$ws.Range("D19").Value = 1.043985076465649
$ws.Range("E19").Value = 1.034071487210028
$ws.Range("F19").Value = 1.054486220627644
$ws.Range("I19").Value = 1.042908857320574
$ws.Range("J19").Value = 1.040972391302472
$ws.Range("K19").Value = 1.047267599877307
$ws.Range("L19").Value = 1.037387730957305
$ws.Range("M19").Value = 1.057733783978134
$ws.Range("N19").Value = 1.042450691499323

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.034593921162754
$ws.Range("D20").Value = 1.043751468175665
$ws.Range("E20").Value = 1.033822059394806
$ws.Range("F20").Value = 1.05419459317904
$ws.Range("I20").Value = 1.04282300093987
$ws.Range("J20").Value = 1.040799508361892
$ws.Range("K20").Value = 1.047098150283389
$ws.Range("L20").Value = 1.037203226172525
$ws.Range("M20").Value = 1.057505776751511
$ws.Range("N20").Value = 1.042277563045138

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.033634492192224
$ws.Range("D21").Value = 1.042992701940217
$ws.Range("E21").Value = 1.033012409080324
$ws.Range("F21").Value = 1.053247729628583
$ws.Range("I21").Value = 1.042542506831388
$ws.Range("J21").Value = 1.040237365403258
$ws.Range("K21").Value = 1.046547008927303
$ws.Range("L21").Value = 1.036603701430914
$ws.Range("M21").Value = 1.056764839390985
$ws.Range("N21").Value = 1.041714621779024

$ws.Range("B22").Value = 1.019999999999999
$ws.Range("C22").Value = 1.033031915572006
$ws.Range("D22").Value = 1.042516194045244
$ws.Range("E22").Value = 1.032504327618912
$ws.Range("F22").Value = 1.052653363478276
$ws.Range("I22").Value = 1.0423651055572
$ws.Range("J22").Value = 1.039883865970997
$ws.Range("K22").Value = 1.046200303840339
$ws.Range("L22").Value = 1.036227007250909
$ws.Range("M22").Value = 1.056299248617305
$ws.Range("N22").Value = 1.041360620337022

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.033351290285389
$ws.Range("D23").Value = 1.042768746366862
$ws.Range("E23").Value = 1.032773577930102
$ws.Range("F23").Value = 1.052968356153508
$ws.Range("I23").Value = 1.042459247974574
$ws.Range("J23").Value = 1.040071267777334
$ws.Range("K23").Value = 1.046384115459063
$ws.Range("L23").Value = 1.036426675738114
$ws.Range("M23").Value = 1.056546041351105
$ws.Range("N23").Value = 1.04154828827542

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.034610017096476
$ws.Range("D24").Value = 1.043764198331382
$ws.Range("E24").Value = 1.033835649733521
$ws.Range("F24").Value = 1.054210483681771
$ws.Range("I24").Value = 1.042827685712641
$ws.Range("J24").Value = 1.040808931691077
$ws.Range("K24").Value = 1.047107387083865
$ws.Range("L24").Value = 1.037213281436695
$ws.Range("M24").Value = 1.057518203054704
$ws.Range("N24").Value = 1.042286999756532

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.036073497008339
$ws.Range("D25").Value = 1.044921738135638
$ws.Range("E25").Value = 1.035072322490757
$ws.Range("F25").Value = 1.055656039169866
$ws.Range("I25").Value = 1.043250666761177
$ws.Range("J25").Value = 1.041664656085761
$ws.Range("K25").Value = 1.047945873652395
$ws.Range("L25").Value = 1.038127140415362
$ws.Range("M25").Value = 1.058647446775342
$ws.Range("N25").Value = 1.043143939378008
